# Powerpoint writer: consolidate text run nodes.
# Merge adjacent "word" + "trailing space" runs into single run nodes
# (without altering the overall visible text), matching the new writer
# behaviour that groups a word with the space that follows it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1 ("Testing custom properties" title) ---------------------
$tr1 = $s.Shapes.Item(1).TextFrame.TextRange

# "Testing" + " " -> "Testing "
$tr1.Characters(1, 8).Text = "Testing "
# "custom" + " " -> "custom "
$tr1.Characters(9, 7).Text = "custom "

# --- Shape 2 ("This is a subtitle" / "A. M." subtitle) ----------------
$tr2 = $s.Shapes.Item(2).TextFrame.TextRange

# "This" + " " -> "This "
$tr2.Characters(1, 5).Text = "This "
# "is" + " " -> "is "
$tr2.Characters(6, 3).Text = "is "
# "a" + " " -> "a "
$tr2.Characters(9, 2).Text = "a "
# "A." + " " -> "A. "
$tr2.Characters(21, 3).Text = "A. "
